$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6454
$ws1.Range("F8").Value = 533
$ws1.Range("F9").Value = 89
$ws1.Range("F10").Value = 77
$ws1.Range("F13").Value = 376
$ws1.Range("F14").Value = 944
$ws1.Range("F15").Value = 3158
$ws1.Range("F17").Value = 192
$ws1.Range("F18").Value = 1829
$ws1.Range("G18").Value = 70
$ws1.Range("F19").Value = 23

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6454
$ws4.Range("F9").Value = 533
$ws4.Range("F10").Value = 89
$ws4.Range("F11").Value = 77
$ws4.Range("F14").Value = 376
$ws4.Range("F15").Value = 944
$ws4.Range("F16").Value = 3158
$ws4.Range("F18").Value = 192
$ws4.Range("F19").Value = 1829
$ws4.Range("G19").Value = 70
$ws4.Range("F20").Value = 23
